$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "Frax" row at position 34 ---------------------------------
# This shifts the old rows 34..50 (HuobiToken .. Aave) down to 35..51, and the
# previous last row (51, "NEARProtocol") is pushed out to row 52.
$ws.Rows("34").Insert()

# The table keeps a fixed extent (A1:E51), so the row that fell off the bottom
# ("NEARProtocol", now at 52) is removed.
$ws.Rows("52").Delete()

# Row-insert leaves the new row's index cell (column A) without the shared
# numeric/border style used by the rest of column A; copy formats from a
# neighbor so A34 matches the look of A2:A51.
$ws.Range("A33").Copy()
$ws.Range("A34").PasteSpecial(-4122)

# Column A holds a simple 0-based running counter; re-number it for the
# inserted row and every row that shifted down by one.
for ($r = 34; $r -le 51; $r++) {
    $ws.Range("A$r").Value = ($r - 2)
}

$ws.Range("D2").Value = "26.482.79"
$ws.Range("D3").Value = "1.724.92"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9971"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.96"
$ws.Range("E5").Value = "  -1.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9977"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4900"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2608"
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06202"
$ws.Range("E9").Value = "  -0.23%  "
$ws.Range("D10").Value = "1.728.71"
$ws.Range("E10").Value = "  -0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06971"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.528"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6019"
$ws.Range("E14").Value = "  -1.06%  "
$ws.Range("E15").Value = "  -0.27%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9974"
$ws.Range("E16").Value = "  -0.27%  "
$ws.Range("D17").Value = "26.476.34"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007178"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D21").Value = "1.953.62"
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.449"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.524"
$ws.Range("E23").Value = "  -2.50%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.110"
$ws.Range("E24").Value = "  -2.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.67"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.26"
$ws.Range("E26").Value = "  -1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.408"
$ws.Range("E27").Value = "  -0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.748"
$ws.Range("E28").Value = "  -1.27%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.60"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.921"
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08009"
$ws.Range("E31").Value = "  -0.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.643"
$ws.Range("E32").Value = "  -1.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04485"
$ws.Range("E33").Value = "  -1.73%  "
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9966"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.595"
$ws.Range("E35").Value = "  -0.76%  "
$ws.Range("E36").Value = "  -0.37%  "
$ws.Range("E37").Value = "  -1.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9396"
$ws.Range("E38").Value = "  +4.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.996"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.387"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9972"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01485"
$ws.Range("E42").Value = "  -1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.85"
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.403"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3849"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.895"
$ws.Range("E46").Value = "  -0.53%  "
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05367"
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "30.58"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.729"
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.48"
$ws.Range("E51").Value = "  -0.26%  "
